$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (C) holds plain text like "$199" (not an Excel currency
# number) in the source data, and the Threads column (F) holds plain-text
# numbers like "16". Force both to Text format before writing so Excel
# stores the literal strings instead of auto-converting them to a currency
# value / a real number.
$ws.Range("C2:C38").NumberFormat = "@"
$ws.Range("F2:F38").NumberFormat = "@"

# Row 2
$ws.Range("C2").Value = '$199'

# Row 3
$ws.Range("C3").Value = '$199'

# Row 4
$ws.Range("C4").Value = '$598'

# Row 5
$ws.Range("C5").Value = '$206'

# Row 6
$ws.Range("C6").Value = '$399'

# Row 7
$ws.Range("C7").Value = '$179'

# Row 8
$ws.Range("C8").Value = '$262'

# Row 9
$ws.Range("C9").Value = '$299'

# Row 10
$ws.Range("C10").Value = '$779'

# Row 11
$ws.Range("C11").Value = '$308'

# Row 12
$ws.Range("C12").Value = '$313'

# Row 13
$ws.Range("C13").Value = '$164'

# Row 14
$ws.Range("C14").Value = '$371'

# Row 15
$ws.Range("C15").Value = '$272'

# Row 16
$ws.Range("C16").Value = '$312'

# Row 17
$ws.Range("C17").Value = '$139'

# Row 18
$ws.Range("C18").Value = '$193'

# Row 19
$ws.Range("C19").Value = '$399'

# Row 20
$ws.Range("C20").Value = '$151'

# Row 21
$ws.Range("C21").Value = '$549'

# Row 22
$ws.Range("C22").Value = '$169'

# Row 23
$ws.Range("C23").Value = '$259'

# Row 24
$ws.Range("C24").Value = '$249'

# Row 25
$ws.Range("C25").Value = '$149'

# Row 26
$ws.Range("C26").Value = '$134'

# Row 27
$ws.Range("C27").Value = '$549'

# Row 28
$ws.Range("C28").Value = '$317'

# Row 29
$ws.Range("B29").Value = 'Core i7-11700KF'
$ws.Range("C29").Value = '$298'
$ws.Range("E29").Value = '8-Core'
$ws.Range("F29").Value = '16'
$ws.Range("H29").Value = '5.0 GHz'

# Row 30
$ws.Range("B30").Value = 'Core i3-10100F'
$ws.Range("C30").Value = '$78'
$ws.Range("D30").Value = 'LGA 1200'
$ws.Range("G30").Value = '3.6 GHz'
$ws.Range("H30").Value = '4.30 GHz'

# Row 31
$ws.Range("B31").Value = 'Core i3-12100F'
$ws.Range("C31").Value = '$107'
$ws.Range("D31").Value = 'LGA 1700'
$ws.Range("E31").Value = 'Quad-Core'
$ws.Range("F31").Value = '8'
$ws.Range("G31").Value = '3.3 GHz'
$ws.Range("H31").Value = '4.3 GHz'

# Row 32
$ws.Range("B32").Value = 'Core i7-10700K'
$ws.Range("C32").Value = '$328'
$ws.Range("D32").Value = 'LGA 1200'
$ws.Range("E32").Value = '8-Core'
$ws.Range("F32").Value = '16'
$ws.Range("G32").Value = '3.8 GHz'
$ws.Range("H32").Value = '5.10 GHz'

# Row 33
$ws.Range("B33").Value = 'Core i7-12700KF'
$ws.Range("C33").Value = '$363'
$ws.Range("D33").Value = 'LGA 1700'
$ws.Range("E33").Value = '12-Core (8P+4E)'
$ws.Range("F33").Value = '20'
$ws.Range("G33").Value = 'P-core Base Frequency: 3.6 GHzE-core Base Frequency: 2.7 GHz'
$ws.Range("H33").Value = 'Intel Turbo Boost Max Technology 3.0 Frequency: Up to 5.0 GHzSingle P-core Turbo Frequency: Up to 4.9 GHzSingle E-core Turbo Frequency: Up to 3.8 GHz'

# Row 34
$ws.Range("B34").Value = 'Core i5-10600K'
$ws.Range("C34").Value = '$203'
$ws.Range("G34").Value = '4.1 GHz'
$ws.Range("H34").Value = '4.80 GHz'

# Row 35
$ws.Range("B35").Value = 'Core i5-11600KF'
$ws.Range("C35").Value = '$209'
$ws.Range("E35").Value = '6-Core'
$ws.Range("F35").Value = '12'
$ws.Range("G35").Value = '3.9 GHz'
$ws.Range("H35").Value = '4.9 GHz'

# Row 36
$ws.Range("B36").Value = 'Core i3-10105'
$ws.Range("C36").Value = '$111'
$ws.Range("E36").Value = 'Quad-Core'
$ws.Range("F36").Value = '8'
$ws.Range("H36").Value = '4.4 GHz'

# Row 37
$ws.Range("B37").Value = 'Core i9-10900K'
$ws.Range("C37").Value = '$389'
$ws.Range("D37").Value = 'LGA 1200'
$ws.Range("E37").Value = '10-Core'
$ws.Range("F37").Value = '20'
$ws.Range("G37").Value = '3.7 GHz'
$ws.Range("H37").Value = '5.30 GHz'

# Row 38
$ws.Range("B38").Value = 'Core i9-12900'
$ws.Range("C38").Value = '$508'
$ws.Range("D38").Value = 'LGA 1700'
$ws.Range("E38").Value = '16-Core (8P+8E)'
$ws.Range("F38").Value = '24'
$ws.Range("G38").Value = 'P-core Base Frequency: 2.4 GHzE-core Base Frequency: 1.8 GHz'
$ws.Range("H38").Value = 'Intel Turbo Boost Max Technology 3.0 Frequency: Up to 5.1 GHzP-core Turbo Frequency: Up to 5.0 GHzE-core Turbo Frequency: Up to 3.8 GHz'
